# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from 2024-04-08 (serial 45390) to 2024-04-09 (serial 45391).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45390) {
        $cell.Value2 = 45391
    }
}
